# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.656.89'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.523.08'
$ws.Range("E3").Value = '  -1.06%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.17'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.49'
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.516.27'
$ws.Range("E7").Value = '  -1.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.609'
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.196'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.24'
$ws.Range("E11").Value = '  +3.50%  '
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.50'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000276'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.096.76'
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.44'
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '617.78'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.526.38'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.719.37'
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("E20").Value = '  +1.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.72'
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.884'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.98'
$ws.Range("E23").Value = '  -6.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.74'
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.54'
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.84'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.17'
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.04'
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -4.39%  '
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.84'
$ws.Range("E34").Value = '  -2.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '621.15'
$ws.Range("E35").Value = '  +8.17%  '
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("B37").Value = 'Cosmos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.84'
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0483'
$ws.Range("E38").Value = '  +2.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.49'
$ws.Range("E39").Value = '  -4.00%  '
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.371.46'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0737'
$ws.Range("E44").Value = '  +4.06%  '
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.94'
$ws.Range("E46").Value = '  -2.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.25'
$ws.Range("E47").Value = '  -2.56%  '
$ws.Range("E48").Value = '  -2.70%  '
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.96'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("E51").Value = '  +0.00%  '
